$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 892, shifting rows
# 892:922 down to 893:923 (dimension grows from A1:R922 to A1:R923).
$ws.Rows("892:892").Insert()

# Populate the newly inserted row 892 with the new "Camote" record.
$ws.Cells.Item(892, 1).Value = 8
$ws.Cells.Item(892, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(892, 3).Value = "Coquimbo"
$ws.Cells.Item(892, 4).Value = 45075
$ws.Cells.Item(892, 5).Value = 4
$ws.Cells.Item(892, 6).Value = 100112045
$ws.Cells.Item(892, 7).Value = "Zapallo"
$ws.Cells.Item(892, 8).Value = "Camote"
$ws.Cells.Item(892, 9).Value = "1a (guarda)"
$ws.Cells.Item(892, 10).Value = 1000
$ws.Cells.Item(892, 11).Value = 400
$ws.Cells.Item(892, 12).Value = 500
$ws.Cells.Item(892, 13).Value = 450
$ws.Cells.Item(892, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(892, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(892, 16).Value = 450
$ws.Cells.Item(892, 17).Value = 1
$ws.Cells.Item(892, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(892, 4).NumberFormat = $ws.Cells.Item(893, 4).NumberFormat
